$wb = $excel.ActiveWorkbook

# The sheet was renamed from "REFERENCES" to "CUSTOMERS" (see commit message:
# "Magic sheet Names moved to ExcelProperties class"). Renaming the sheet via
# the object model also updates every reference to the old name, including
# the hidden _xlnm._FilterDatabase defined name that anchors the autofilter.
$ws = $wb.ActiveSheet
$ws.Name = "CUSTOMERS"
